# Update the "想去人数" (F column) counts across the sheets to the newly
# scraped values, as published to gh-pages at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) -- column F, rows 2-48
$ws1 = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 99
    3  = 1037
    4  = 9234
    5  = 193
    6  = 65
    7  = 1956
    8  = 6380
    9  = 618
    10 = 66
    11 = 9437
    12 = 11046
    13 = 1230
    14 = 1137
    15 = 4908
    16 = 791
    17 = 442
    19 = 327
    20 = 175
    21 = 1330
    22 = 236
    23 = 875
    24 = 1228
    25 = 853
    26 = 3
    27 = 2018
    28 = 422
    29 = 613
    30 = 2647
    31 = 306
    32 = 182
    33 = 1731
    34 = 92
    35 = 1332
    36 = 442
    37 = 42
    38 = 911
    39 = 589
    40 = 15
    41 = 3296
    42 = 234
    44 = 505
    45 = 574
    48 = 235
}
foreach ($row in $exhibitionUpdates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "演出" (Performances) -- column F, row 19
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(19, 6).Value = 9

# Sheet "本地生活" (Local Life) -- column F, row 2
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 5881

# Sheet "全部类型" (All Types) -- column F, rows 2-48
$ws4 = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    2  = 1037
    3  = 9234
    4  = 65
    7  = 6380
    8  = 618
    9  = 9437
    10 = 11046
    12 = 1230
    13 = 1137
    14 = 4908
    15 = 791
    16 = 442
    18 = 175
    20 = 1330
    21 = 236
    22 = 875
    23 = 1228
    24 = 853
    26 = 2018
    27 = 422
    28 = 613
    29 = 2647
    30 = 182
    31 = 1731
    33 = 442
    37 = 911
    40 = 589
    42 = 234
    44 = 505
    45 = 574
    48 = 235
}
foreach ($row in $allTypesUpdates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
